# Add a new "UK" test-data worksheet, modelled on the existing "Poland"
# sheet (same layout/column widths/cell styles), inserted right after it.
# Matches commit "Added Test Data for UK Market".

$wb = $excel.ActiveWorkbook

$poland = $wb.Worksheets.Item("Poland")

# Duplicate "Poland" (keeps column widths, merged cells and cell styles
# identical, and places the copy immediately after it) then rename it.
$poland.Copy([System.Reflection.Missing]::Value, $poland)
$ukSheet = $wb.Worksheets.Item($poland.Index + 1)
$ukSheet.Name = "UK"

# The UK repeaters list has two extra models (P32AR / P32DR) that Poland's
# list doesn't have, listed just before "MZXDR240" -> insert two rows there.
$ukSheet.Range("A15:A16").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown) | Out-Null

# Give the freshly inserted rows the same formatting (border/fill/font) as
# the rest of the list, by copying an existing list cell's format onto them.
$ukSheet.Range("A14").Copy() | Out-Null
$ukSheet.Range("A15:A16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ukSheet.Range("A15").Value = "P32AR"
$ukSheet.Range("A16").Value = "P32DR"

# Market-specific values (order matters only for shared-string ordering).
$ukSheet.Range("B4").Value = "NGC-2741/T3349"
$ukSheet.Range("B2").Value = "UK Market"

# Match the selection left on the new sheet by the author.
$ukSheet.Range("A16").Select() | Out-Null

Write-Output "Added UK worksheet"
